$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.745.71"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.534.50"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.244"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.24"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0578"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "1.757.05"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "1.543.14"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.505"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "26.738.33"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "211.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "1.359.61"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.796"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "1.670.85"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("D50").Value = "0.0₇0974"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0947"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.68%  "
